# Updates "想去人数" (interest count) figures across the workbook.
# Mirrors the upstream "Update gh-pages to output generated at 456a3b4" commit.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 614
$ws1.Range("F6").Value = 414
$ws1.Range("F8").Value = 190
$ws1.Range("F10").Value = 255
$ws1.Range("F11").Value = 6993
$ws1.Range("F12").Value = 72
$ws1.Range("F20").Value = 727
$ws1.Range("F24").Value = 338
$ws1.Range("F27").Value = 30
$ws1.Range("F28").Value = 1972
$ws1.Range("F29").Value = 544
$ws1.Range("F31").Value = 540

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 280

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 319

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 319
$ws4.Range("F3").Value = 614
$ws4.Range("F8").Value = 414
$ws4.Range("F10").Value = 190
$ws4.Range("F12").Value = 255
$ws4.Range("F13").Value = 6993
$ws4.Range("F14").Value = 72
$ws4.Range("F24").Value = 280
$ws4.Range("F27").Value = 727
$ws4.Range("F34").Value = 338
$ws4.Range("F37").Value = 30
$ws4.Range("F38").Value = 1972
$ws4.Range("F39").Value = 544
$ws4.Range("F41").Value = 540
